$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 282, shifting the existing rows (282:300) down to (283:301).
$ws.Rows.Item(282).Insert()

# Populate the newly inserted row 282 with the new weekly price record.
$ws.Cells.Item(282, 1).Value = 10
$ws.Cells.Item(282, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(282, 3).Value = "La Araucanía"
$ws.Cells.Item(282, 4).Value = 44746
$ws.Cells.Item(282, 5).Value = 9
$ws.Cells.Item(282, 6).Value = 100112001
$ws.Cells.Item(282, 7).Value = "Berenjena"
$ws.Cells.Item(282, 8).Value = "Sin especificar"
$ws.Cells.Item(282, 9).Value = "Primera"
$ws.Cells.Item(282, 10).Value = 65
$ws.Cells.Item(282, 11).Value = 12000
$ws.Cells.Item(282, 12).Value = 12000
$ws.Cells.Item(282, 13).Value = 12000
$ws.Cells.Item(282, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(282, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(282, 16).Value = 200
$ws.Cells.Item(282, 17).Value = 60
$ws.Cells.Item(282, 18).Value = "Hortaliza"
